$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting "Model"/"All features"... into column B,
# and the accuracy values into column C.
$ws.Columns("A").Insert()

# Add a new row for the "Linear SVC" classifier results (written first so its
# shared string is registered before "SGD", matching the target string order)
$ws.Range("A6").Value2 = "Linear SVC"
$ws.Range("B6").Value2 = "All features"
$ws.Range("C6").Value2 = 0.8167
$ws.Range("C6").NumberFormat = "0.00%"

# Fill the new column A with the classifier name for each existing data row (2-5): "SGD"
$ws.Range("A2").Value2 = "SGD"
$ws.Range("A3").Value2 = "SGD"
$ws.Range("A4").Value2 = "SGD"
$ws.Range("A5").Value2 = "SGD"

# Match the final selected cell recorded in the saved workbook
$ws.Range("E18").Select() | Out-Null

Write-Host "Done"
